$wb = $excel.ActiveWorkbook

# The sheets "展览" (Exhibitions) and "全部类型" (All types) hold identical
# event-listing tables. The edit:
#   - removes the oldest event (2024.02.19 安徽马娘only) that used to sit in row 2
#   - shifts every later event up by one row
#   - tweaks a handful of numeric "want to go" counts and one cover image URL
#   - drops the now-duplicated trailing row, shrinking the table from 7 to 6 rows
#
# Rather than literally deleting row 2 (which, via COM, would also renumber the
# manual index values in column A), we overwrite rows 2-6 with their final
# target contents directly, leaving column A untouched, and then delete the
# now-redundant last row (row 7) to shrink the used range.
# (Note: this runtime's PowerShell only binds positional parameters, so the
# helper below avoids named "-param value" syntax.)

function Set-EventRow($ws, $row, $date, $title, $place, $timeRange, $want, $minPrice, $link, $cover) {
    # Column B holds plain "YYYY.MM.DD" text. Typed verbatim, Excel's COM
    # layer auto-recognizes that pattern as a date and silently converts it
    # to a date serial number. Prefixing with an apostrophe (exactly what a
    # user would do in the Excel UI) forces it to stay literal text.
    $ws.Range("B" + $row).Value = "'" + $date
    $ws.Range("C" + $row).Value = $title
    $ws.Range("D" + $row).Value = $place
    $ws.Range("E" + $row).Value = $timeRange
    $ws.Range("F" + $row).Value = $want
    $ws.Range("G" + $row).Value = $minPrice
    $ws.Range("H" + $row).Value = $link
    $ws.Range("I" + $row).Value = $cover
}

function Update-EventSheet($ws) {
    Set-EventRow $ws 2 "2024.03.02" "合肥·星芒1.5动漫嘉年华" `
        "山西路与太原路交叉口 挥动体育" `
        "2024.03.02 09:30-03.02 17:30" `
        1279 49.5 `
        "https://show.bilibili.com/platform/detail.html?id=81267" `
        "//i0.hdslb.com/bfs/openplatform/202401/GWidiefU1706003134747.jpeg"

    Set-EventRow $ws 3 "2024.03.16" "合肥·CW国潮动漫游戏嘉年华" `
        "南京路与庐州大道交汇处 合肥滨湖国际会展中心" `
        "2024.03.16 09:30-03.17 17:00" `
        1600 65 `
        "https://show.bilibili.com/platform/detail.html?id=81284" `
        "//i0.hdslb.com/bfs/openplatform/202401/38B92fWF1705995243803.jpeg"

    Set-EventRow $ws 4 "2024.03.23" "合肥·原&铁&崩 only展" `
        "金寨路与天堂窄路交叉口 梵木艺术中心" `
        "2024.03.23 09:00-03.23 17:00" `
        63 58 `
        "https://show.bilibili.com/platform/detail.html?id=81574" `
        "//i2.hdslb.com/bfs/openplatform/202401/0V5uyX6C1706697212904.png"

    Set-EventRow $ws 5 "2024.04.04" "合肥· 第二届漫画城市动漫展 -故事再次开始" `
        "凤淮路与固镇路西北角 庐阳全民健身中心" `
        "2024.04.04 09:00-04.05 17:00" `
        6199 60 `
        "https://show.bilibili.com/platform/detail.html?id=78898" `
        "//i2.hdslb.com/bfs/openplatform/202402/3NgyB9761708333056023.jpeg"

    Set-EventRow $ws 6 "2024.05.18" "合肥·梦时空SPO1动漫展" `
        "阜阳路16号 银瑞林国际大酒店" `
        "2024.05.18 10:00-05.18 17:00" `
        105 60 `
        "https://show.bilibili.com/platform/detail.html?id=80207" `
        "//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg"

    # Row 7 (old "梦时空SPO1" row, now duplicated into row 6 above) is no
    # longer needed; removing it shrinks the sheet's dimension from I7 to I6.
    $ws.Rows.Item(7).Delete()
}

Update-EventSheet $wb.Worksheets.Item("展览")
Update-EventSheet $wb.Worksheets.Item("全部类型")
